# Add a new bullet to the "Problems" list, right after the paragraph that
# talks about expanding the dataset with a moving window.

$d = $word.ActiveDocument

$anchorText = "Expanded our dataset by using a moving window of 100 transactions for the most recent 1000 transactions."

# Locate the anchor paragraph via Find so the script isn't dependent on
# absolute paragraph indices.
$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false,
                                    $false, $false, $true, 1, $false,
                                    "", 0)

if ($found) {
    $anchorIndex = $searchRange.Paragraphs.Last.Index
} else {
    # Fallback: just use the last paragraph of the document.
    $anchorIndex = $d.Paragraphs.Count
}

# Re-fetch the paragraph from the document's own collection (the range
# returned by Find is narrowed to the match, so go back through
# $d.Paragraphs to get the full paragraph object).
$anchorParagraph = $d.Paragraphs.Item($anchorIndex)

# Insert a new paragraph right after it; Word carries over the source
# paragraph's formatting (ListParagraph style + the same numbered-list
# properties), matching the other bullets in this list.
$anchorParagraph.Range.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Item($anchorIndex + 1)
$newParagraph.Range.Text = "Due to the low amounts of data, the random state of the split can effect the overall accuracy by around 10%."
